$d = $word.ActiveDocument

$replacements = @(
    @("949÷2=", "589÷5="),
    @("386÷5=", "647÷7="),
    @("114÷4=", "636÷6="),
    @("554÷8=", "310÷8="),
    @("886÷9=", "338÷9="),
    @("465÷9=", "692÷6="),
    @("668÷8=", "423÷7="),
    @("260÷9=", "259÷7="),
    @("355÷3=", "415÷4="),
    @("752÷2=", "369÷9="),
    @("779÷3=", "141÷8="),
    @("257÷8=", "597÷7="),
    @("684÷8=", "587÷5="),
    @("165÷7=", "703÷9="),
    @("214÷2=", "114÷3="),
    @("831÷4=", "299÷9="),
    @("936÷9=", "420÷9="),
    @("357÷4=", "544÷9="),
    @("453÷4=", "432÷6="),
    @("155÷2=", "474÷6="),
    @("461÷2=", "708÷7="),
    @("224÷8=", "914÷9="),
    @("842÷3=", "827÷2="),
    @("408÷6=", "453÷3="),
    @("184÷7=", "137÷9=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
